# Auto-generated Excel COM-interop edit script
# Updates crypto market data across all three worksheets to match the
# 2024-11-22 04:36:14 data refresh described in the commit message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Top 50 Cryptocurrencies"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

# Row 2
$ws1.Cells.Item(2, 3).Value = 98971
$ws1.Cells.Item(2, 4).Value = 1958688961585
$ws1.Cells.Item(2, 5).Value = 90706112783
$ws1.Cells.Item(2, 6).Value = 2.72787

# Row 3
$ws1.Cells.Item(3, 3).Value = 3414.69
$ws1.Cells.Item(3, 4).Value = 409915509517
$ws1.Cells.Item(3, 5).Value = 58463062946
$ws1.Cells.Item(3, 6).Value = 9.616110000000001

# Row 4
$ws1.Cells.Item(4, 3).Value = 1.001
$ws1.Cells.Item(4, 4).Value = 130869523334
$ws1.Cells.Item(4, 5).Value = 144923567712
$ws1.Cells.Item(4, 6).Value = 0.16594

# Row 5
$ws1.Cells.Item(5, 3).Value = 261.97
$ws1.Cells.Item(5, 4).Value = 124533585904
$ws1.Cells.Item(5, 5).Value = 14798492535
$ws1.Cells.Item(5, 6).Value = 10.07209

# Row 6
$ws1.Cells.Item(6, 3).Value = 635.71
$ws1.Cells.Item(6, 4).Value = 92659855271
$ws1.Cells.Item(6, 5).Value = 2505513002
$ws1.Cells.Item(6, 6).Value = 5.36418

# Row 7
$ws1.Cells.Item(7, 3).Value = 1.4
$ws1.Cells.Item(7, 4).Value = 80057199926
$ws1.Cells.Item(7, 5).Value = 17726236905
$ws1.Cells.Item(7, 6).Value = 27.44776

# Row 8
$ws1.Cells.Item(8, 3).Value = 0.395193
$ws1.Cells.Item(8, 4).Value = 58107129028
$ws1.Cells.Item(8, 5).Value = 10365661625
$ws1.Cells.Item(8, 6).Value = 3.48193

# Row 9
$ws1.Cells.Item(9, 3).Value = 0.999896
$ws1.Cells.Item(9, 4).Value = 38235535434
$ws1.Cells.Item(9, 5).Value = 14818752221
$ws1.Cells.Item(9, 6).Value = 0.09767000000000001

# Row 10
$ws1.Cells.Item(10, 3).Value = 3413.2
$ws1.Cells.Item(10, 4).Value = 33428135997
$ws1.Cells.Item(10, 5).Value = 148487025
$ws1.Cells.Item(10, 6).Value = 9.59801

# Row 11
$ws1.Cells.Item(11, 3).Value = 0.890408
$ws1.Cells.Item(11, 4).Value = 31913200608
$ws1.Cells.Item(11, 5).Value = 3611716130
$ws1.Cells.Item(11, 6).Value = 14.09951

# Row 12
$ws1.Cells.Item(12, 3).Value = 0.20109
$ws1.Cells.Item(12, 4).Value = 17365839559
$ws1.Cells.Item(12, 5).Value = 1101570407
$ws1.Cells.Item(12, 6).Value = 2.44591

# Row 13
$ws1.Cells.Item(13, 3).Value = 36.47
$ws1.Cells.Item(13, 4).Value = 14930904597
$ws1.Cells.Item(13, 5).Value = 1055426491
$ws1.Cells.Item(13, 6).Value = 8.78505

# Row 14
$ws1.Cells.Item(14, 3).Value = 0.00002507
$ws1.Cells.Item(14, 4).Value = 14768662928
$ws1.Cells.Item(14, 5).Value = 1622090425
$ws1.Cells.Item(14, 6).Value = 5.69772

# Row 15
$ws1.Cells.Item(15, 3).Value = 98944
$ws1.Cells.Item(15, 4).Value = 14459289342
$ws1.Cells.Item(15, 5).Value = 903720736
$ws1.Cells.Item(15, 6).Value = 3.29476

# Row 16
$ws1.Cells.Item(16, 3).Value = 4003.17
$ws1.Cells.Item(16, 4).Value = 14446765079
$ws1.Cells.Item(16, 5).Value = 167820354
$ws1.Cells.Item(16, 6).Value = 9.105029999999999

# Row 17
$ws1.Cells.Item(17, 3).Value = 5.57
$ws1.Cells.Item(17, 4).Value = 14197408302
$ws1.Cells.Item(17, 5).Value = 626468405
$ws1.Cells.Item(17, 6).Value = 4.87141

# Row 18
$ws1.Cells.Item(18, 3).Value = 3.64
$ws1.Cells.Item(18, 4).Value = 10370586662
$ws1.Cells.Item(18, 5).Value = 2130512817
$ws1.Cells.Item(18, 6).Value = 3.49773

# Row 19
$ws1.Cells.Item(19, 1).Value = 'WETH'
$ws1.Cells.Item(19, 2).Value = 'weth'
$ws1.Cells.Item(19, 3).Value = 3411.88
$ws1.Cells.Item(19, 4).Value = 9816635144
$ws1.Cells.Item(19, 5).Value = 2314023127
$ws1.Cells.Item(19, 6).Value = 9.716889999999999

# Row 20
$ws1.Cells.Item(20, 1).Value = 'Bitcoin Cash'
$ws1.Cells.Item(20, 2).Value = 'bch'
$ws1.Cells.Item(20, 3).Value = 494.21
$ws1.Cells.Item(20, 4).Value = 9777189234
$ws1.Cells.Item(20, 5).Value = 2240429444
$ws1.Cells.Item(20, 6).Value = 5.43251

# Row 21
$ws1.Cells.Item(21, 3).Value = 15.26
$ws1.Cells.Item(21, 4).Value = 9535373196
$ws1.Cells.Item(21, 5).Value = 1235242963
$ws1.Cells.Item(21, 6).Value = 6.50221

# Row 22
$ws1.Cells.Item(22, 3).Value = 0.00002171
$ws1.Cells.Item(22, 4).Value = 9133951966
$ws1.Cells.Item(22, 5).Value = 7119159661
$ws1.Cells.Item(22, 6).Value = 14.42217

# Row 23
$ws1.Cells.Item(23, 3).Value = 6.24
$ws1.Cells.Item(23, 4).Value = 8987217277
$ws1.Cells.Item(23, 5).Value = 821172786
$ws1.Cells.Item(23, 6).Value = 11.20861

# Row 24
$ws1.Cells.Item(24, 3).Value = 0.285872
$ws1.Cells.Item(24, 4).Value = 8586479712
$ws1.Cells.Item(24, 5).Value = 2347001161
$ws1.Cells.Item(24, 6).Value = 21.45542

# Row 25
$ws1.Cells.Item(25, 4).Value = 8073392142
$ws1.Cells.Item(25, 5).Value = 3453877
$ws1.Cells.Item(25, 6).Value = 2.58685

# Row 26
$ws1.Cells.Item(26, 3).Value = 5.85
$ws1.Cells.Item(26, 4).Value = 7116535637
$ws1.Cells.Item(26, 5).Value = 1022382381
$ws1.Cells.Item(26, 6).Value = 6.4557

# Row 27
$ws1.Cells.Item(27, 3).Value = 90.47
$ws1.Cells.Item(27, 4).Value = 6802702646
$ws1.Cells.Item(27, 5).Value = 1425441937
$ws1.Cells.Item(27, 6).Value = 6.36019

# Row 28
$ws1.Cells.Item(28, 3).Value = 12.14
$ws1.Cells.Item(28, 4).Value = 6472704885
$ws1.Cells.Item(28, 5).Value = 897409636
$ws1.Cells.Item(28, 6).Value = 5.02701

# Row 29
$ws1.Cells.Item(29, 3).Value = 3591.76
$ws1.Cells.Item(29, 4).Value = 6158488048
$ws1.Cells.Item(29, 5).Value = 103343723
$ws1.Cells.Item(29, 6).Value = 10.11175

# Row 30
$ws1.Cells.Item(30, 3).Value = 9.49
$ws1.Cells.Item(30, 4).Value = 5685341830
$ws1.Cells.Item(30, 5).Value = 854581222
$ws1.Cells.Item(30, 6).Value = 8.99269

# Row 31
$ws1.Cells.Item(31, 3).Value = 0.195381
$ws1.Cells.Item(31, 4).Value = 5297070073
$ws1.Cells.Item(31, 5).Value = 116770671
$ws1.Cells.Item(31, 6).Value = 10.9334

# Row 32
$ws1.Cells.Item(32, 3).Value = 0.997246
$ws1.Cells.Item(32, 4).Value = 5236377315
$ws1.Cells.Item(32, 5).Value = 91094
$ws1.Cells.Item(32, 6).Value = -0.13942

# Row 33
$ws1.Cells.Item(33, 3).Value = 0.132292
$ws1.Cells.Item(33, 4).Value = 5045007831
$ws1.Cells.Item(33, 5).Value = 870484679
$ws1.Cells.Item(33, 6).Value = 5.95544

# Row 34
$ws1.Cells.Item(34, 4).Value = 4601900585
$ws1.Cells.Item(34, 5).Value = 269764293
$ws1.Cells.Item(34, 6).Value = 8.44617

# Row 35
$ws1.Cells.Item(35, 3).Value = 28.28
$ws1.Cells.Item(35, 4).Value = 4216363242
$ws1.Cells.Item(35, 5).Value = 905292354
$ws1.Cells.Item(35, 6).Value = 8.049530000000001

# Row 36
$ws1.Cells.Item(36, 3).Value = 0.00005229
$ws1.Cells.Item(36, 4).Value = 3926993569
$ws1.Cells.Item(36, 5).Value = 1801299053
$ws1.Cells.Item(36, 6).Value = 6.7522

# Row 37
$ws1.Cells.Item(37, 3).Value = 7.45
$ws1.Cells.Item(37, 4).Value = 3853285644
$ws1.Cells.Item(37, 5).Value = 447297571
$ws1.Cells.Item(37, 6).Value = 2.48477

# Row 38
$ws1.Cells.Item(38, 3).Value = 0.151479
$ws1.Cells.Item(38, 4).Value = 3820698917
$ws1.Cells.Item(38, 5).Value = 155603978
$ws1.Cells.Item(38, 6).Value = 1.39954

# Row 39
$ws1.Cells.Item(39, 3).Value = 513.33
$ws1.Cells.Item(39, 4).Value = 3789399208
$ws1.Cells.Item(39, 5).Value = 287599186
$ws1.Cells.Item(39, 6).Value = 6.01779

# Row 40
$ws1.Cells.Item(40, 3).Value = 0.475024
$ws1.Cells.Item(40, 4).Value = 3775430360
$ws1.Cells.Item(40, 5).Value = 483340643
$ws1.Cells.Item(40, 6).Value = 9.5351

# Row 41
$ws1.Cells.Item(41, 3).Value = 1.002
$ws1.Cells.Item(41, 4).Value = 3687437906
$ws1.Cells.Item(41, 5).Value = 230417226
$ws1.Cells.Item(41, 6).Value = 0.1649

# Row 42
$ws1.Cells.Item(42, 3).Value = 24.79
$ws1.Cells.Item(42, 4).Value = 3573536978
$ws1.Cells.Item(42, 5).Value = 42132678
$ws1.Cells.Item(42, 6).Value = 3.22954

# Row 43
$ws1.Cells.Item(43, 3).Value = 3.84
$ws1.Cells.Item(43, 4).Value = 3470901788
$ws1.Cells.Item(43, 5).Value = 307312059
$ws1.Cells.Item(43, 6).Value = 6.63508

# Row 44
$ws1.Cells.Item(44, 3).Value = 1
$ws1.Cells.Item(44, 4).Value = 3443844748
$ws1.Cells.Item(44, 5).Value = 161352064
$ws1.Cells.Item(44, 6).Value = 0.29964

# Row 45
$ws1.Cells.Item(45, 3).Value = 3.41
$ws1.Cells.Item(45, 4).Value = 3403124441
$ws1.Cells.Item(45, 5).Value = 1289304237
$ws1.Cells.Item(45, 6).Value = 9.82887

# Row 46
$ws1.Cells.Item(46, 4).Value = 3362395368
$ws1.Cells.Item(46, 5).Value = 497643432
$ws1.Cells.Item(46, 6).Value = 5.78979

# Row 47
$ws1.Cells.Item(47, 3).Value = 0.79455
$ws1.Cells.Item(47, 4).Value = 3247225217
$ws1.Cells.Item(47, 5).Value = 1700971452
$ws1.Cells.Item(47, 6).Value = 15.00628

# Row 48
$ws1.Cells.Item(48, 3).Value = 160.49
$ws1.Cells.Item(48, 4).Value = 2959603040
$ws1.Cells.Item(48, 5).Value = 83651753
$ws1.Cells.Item(48, 6).Value = -0.51831

# Row 49
$ws1.Cells.Item(49, 3).Value = 1.96
$ws1.Cells.Item(49, 4).Value = 2947490991
$ws1.Cells.Item(49, 5).Value = 398272260
$ws1.Cells.Item(49, 6).Value = 4.18586

# Row 50
$ws1.Cells.Item(50, 3).Value = 4.71
$ws1.Cells.Item(50, 4).Value = 2827521355
$ws1.Cells.Item(50, 5).Value = 590206830
$ws1.Cells.Item(50, 6).Value = 9.32207

# Row 51
$ws1.Cells.Item(51, 3).Value = 46.75
$ws1.Cells.Item(51, 4).Value = 2799927722
$ws1.Cells.Item(51, 5).Value = 20401183
$ws1.Cells.Item(51, 6).Value = 6.80962

# ---------------------------------------------------------------------------
# Sheet 2: "Top 5 by Market Cap"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")

$ws2.Cells.Item(2, 2).Value = 1958688961585
$ws2.Cells.Item(3, 2).Value = 409915509517
$ws2.Cells.Item(4, 2).Value = 130869523334
$ws2.Cells.Item(5, 2).Value = 124533585904
$ws2.Cells.Item(6, 2).Value = 92659855271

# ---------------------------------------------------------------------------
# Sheet 3: "Summary"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Summary")

$ws3.Cells.Item(2, 2).NumberFormat = "@"
$ws3.Cells.Item(2, 2).Value = '$4363.03'
$ws3.Cells.Item(3, 2).Value = 'XRP (27.45%)'
$ws3.Cells.Item(4, 2).Value = 'Monero (-0.52%)'
